$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 219, shifting existing rows 219:232 down to 220:233.
$ws.Rows("219:219").Insert()

# Populate the new row 219 with the new price-report record.
$ws.Cells.Item(219, 1).Value = 10
$ws.Cells.Item(219, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(219, 3).Value = "La Araucanía"
$ws.Cells.Item(219, 4).Value = 44931
$ws.Cells.Item(219, 5).Value = 9
$ws.Cells.Item(219, 6).Value = 100114007
$ws.Cells.Item(219, 7).Value = "Jengibre"
$ws.Cells.Item(219, 8).Value = "Sin especificar"
$ws.Cells.Item(219, 9).Value = "Primera"
$ws.Cells.Item(219, 10).Value = 150
$ws.Cells.Item(219, 11).Value = 20000
$ws.Cells.Item(219, 12).Value = 20000
$ws.Cells.Item(219, 13).Value = 20000
$ws.Cells.Item(219, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(219, 15).Value = "Perú"
$ws.Cells.Item(219, 16).Value = 1538
$ws.Cells.Item(219, 17).Value = 13
$ws.Cells.Item(219, 18).Value = "Hortaliza"

# Match the date style used by the other rows in column D.
$ws.Cells.Item(219, 4).NumberFormat = $ws.Cells.Item(220, 4).NumberFormat
